$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 226
$ws.Range("F4").Value = 810
$ws.Range("F6").Value = 401
$ws.Range("F7").Value = 555
$ws.Range("F8").Value = 212
$ws.Range("F10").Value = 333
$ws.Range("F11").Value = 126
$ws.Range("F12").Value = 603
$ws.Range("F14").Value = 1756
$ws.Range("F15").Value = 320
$ws.Range("F16").Value = 2273
$ws.Range("F17").Value = 289
$ws.Range("F18").Value = 485
$ws.Range("F19").Value = 40

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 220
$ws.Range("F7").Value = 471
$ws.Range("F13").Value = 83

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5279
$ws.Range("F4").Value = 189

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5279
$ws.Range("F6").Value = 189
$ws.Range("F7").Value = 226
$ws.Range("F8").Value = 220
$ws.Range("F12").Value = 471
$ws.Range("F13").Value = 810
$ws.Range("F17").Value = 401
$ws.Range("F18").Value = 555
$ws.Range("F19").Value = 212
$ws.Range("F22").Value = 333
$ws.Range("F23").Value = 126
$ws.Range("F26").Value = 603
$ws.Range("F28").Value = 83
$ws.Range("F29").Value = 1756
$ws.Range("F30").Value = 320
$ws.Range("F31").Value = 2273
$ws.Range("F33").Value = 289
$ws.Range("F34").Value = 485
$ws.Range("F35").Value = 40
